$d = $word.ActiveDocument

$replacements = @(
    @{old="156÷6="; new="390÷5="},
    @{old="631÷2="; new="206÷9="},
    @{old="923÷7="; new="689÷3="},
    @{old="659÷4="; new="736÷4="},
    @{old="405÷2="; new="595÷8="},
    @{old="645÷2="; new="103÷4="},
    @{old="614÷5="; new="621÷7="},
    @{old="452÷3="; new="639÷9="},
    @{old="292÷6="; new="173÷5="},
    @{old="693÷8="; new="704÷3="},
    @{old="253÷8="; new="512÷2="},
    @{old="423÷4="; new="833÷5="},
    @{old="211÷3="; new="733÷8="},
    @{old="821÷9="; new="324÷2="},
    @{old="911÷5="; new="386÷5="},
    @{old="428÷2="; new="742÷4="},
    @{old="387÷2="; new="901÷2="},
    @{old="885÷5="; new="926÷9="},
    @{old="836÷2="; new="477÷9="},
    @{old="443÷8="; new="296÷9="},
    @{old="800÷5="; new="301÷5="},
    @{old="985÷5="; new="187÷5="},
    @{old="699÷7="; new="482÷7="},
    @{old="332÷8="; new="142÷9="},
    @{old="196÷2="; new="696÷4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
